# Update "想去人数" (F) and, where applicable, "最低票价" (G) figures across the
# 展览 / 演出 / 本地生活 / 全部类型 sheets, matching a refreshed scrape snapshot.

$wb = $excel.ActiveWorkbook

function Set-FValue {
    param($ws, [int]$row, [double]$value)
    $ws.Cells.Item($row, 6).Value = $value
}

function Set-GValue {
    # G holds the lowest-ticket-price column, stored as text even when the
    # label is numeric-looking (e.g. "65"). Force text entry via a
    # temporary "@" number format, then restore the cell's original
    # (unstyled) look so only the value actually changes.
    param($ws, [int]$row, [string]$value)
    $cell = $ws.Cells.Item($row, 7)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---- 展览 (sheet 1) ----
$ws1 = $wb.Worksheets.Item("展览")

Set-FValue $ws1 2  1286
Set-FValue $ws1 3  2090
Set-FValue $ws1 4  431
Set-FValue $ws1 6  433
Set-FValue $ws1 8  526
Set-FValue $ws1 9  142
Set-FValue $ws1 11 173
Set-FValue $ws1 12 811
Set-FValue $ws1 13 58
Set-FValue $ws1 15 4350
Set-FValue $ws1 17 855
Set-FValue $ws1 18 632
Set-FValue $ws1 20 728
Set-FValue $ws1 21 1443
Set-FValue $ws1 22 51
Set-GValue $ws1 22 "65"
Set-FValue $ws1 23 682
Set-FValue $ws1 26 217

# ---- 演出 (sheet 2) ----
$ws2 = $wb.Worksheets.Item("演出")
Set-FValue $ws2 2 41

# ---- 本地生活 (sheet 3) ----
$ws3 = $wb.Worksheets.Item("本地生活")
Set-FValue $ws3 2 122

# ---- 全部类型 (sheet 4) ----
$ws4 = $wb.Worksheets.Item("全部类型")

Set-FValue $ws4 2  122
Set-FValue $ws4 3  41
Set-FValue $ws4 5  1286
Set-FValue $ws4 6  2090
Set-FValue $ws4 7  431
Set-FValue $ws4 9  433
Set-FValue $ws4 11 526
Set-FValue $ws4 12 142
Set-FValue $ws4 14 173
Set-FValue $ws4 15 811
Set-FValue $ws4 16 58
Set-FValue $ws4 20 4350
Set-FValue $ws4 22 855
Set-FValue $ws4 23 632
Set-FValue $ws4 25 728
Set-FValue $ws4 26 1443
Set-FValue $ws4 27 51
Set-GValue $ws4 27 "65"
Set-FValue $ws4 28 682
Set-FValue $ws4 31 217
